$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp (row 1)
$ws.Cells.Item(1,1).Value = "Datos actualizados a 9 de Septiembre de 2020 a las 23:27"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4,2).Value = 6543635
$ws.Cells.Item(4,3).Value = 29404
$ws.Cells.Item(4,4).Value = 3827529
$ws.Cells.Item(4,5).Value = 2521161
$ws.Cells.Item(4,7).Value = 915
$ws.Cells.Item(4,8).Value = 194945

# Row 6: Brasil -> Brasil
$ws.Cells.Item(6,2).Value = 4197889
$ws.Cells.Item(6,3).Value = 32765
$ws.Cells.Item(6,5).Value = 672116
$ws.Cells.Item(6,7).Value = 1022
$ws.Cells.Item(6,8).Value = 128539

# Row 17: Francia -> Francia
$ws.Cells.Item(17,4).Value = 88524
$ws.Cells.Item(17,5).Value = 224783

# Row 29: Canada -> Canada
$ws.Cells.Item(29,2).Value = 134096
$ws.Cells.Item(29,3).Value = 348
$ws.Cells.Item(29,5).Value = 6997

# Row 43: Guatemala -> Guatemala
$ws.Cells.Item(43,2).Value = 79622
$ws.Cells.Item(43,3).Value = 901
$ws.Cells.Item(43,4).Value = 68308
$ws.Cells.Item(43,5).Value = 8417
$ws.Cells.Item(43,7).Value = 7
$ws.Cells.Item(43,8).Value = 2897

# Row 53: Singapur -> Barein
$ws.Cells.Item(53,1).Value = "Barein"
$ws.Cells.Item(53,2).Value = 57450
$ws.Cells.Item(53,3).Value = 672
$ws.Cells.Item(53,4).Value = 51819
$ws.Cells.Item(53,5).Value = 5428
$ws.Cells.Item(53,7).Value = 1
$ws.Cells.Item(53,8).Value = 203

# Row 54: Barein -> Singapur
$ws.Cells.Item(54,1).Value = "Singapur"
$ws.Cells.Item(54,2).Value = 57166
$ws.Cells.Item(54,3).Value = 75
$ws.Cells.Item(54,4).Value = 56492
$ws.Cells.Item(54,5).Value = 647
$ws.Cells.Item(54,7).Value = 0
$ws.Cells.Item(54,8).Value = 27

# Row 57: Costa Rica -> Costa Rica
$ws.Cells.Item(57,2).Value = 51224
$ws.Cells.Item(57,3).Value = 1327
$ws.Cells.Item(57,4).Value = 19544
$ws.Cells.Item(57,5).Value = 31137
$ws.Cells.Item(57,7).Value = 12
$ws.Cells.Item(57,8).Value = 543

# Row 82: Costa de Marfil -> Costa de Marfil
$ws.Cells.Item(82,2).Value = 18815
$ws.Cells.Item(82,3).Value = 37
$ws.Cells.Item(82,4).Value = 17770
$ws.Cells.Item(82,5).Value = 926

# Row 95: Guinea -> Guinea
$ws.Cells.Item(95,2).Value = 9885
$ws.Cells.Item(95,3).Value = 37
$ws.Cells.Item(95,4).Value = 9068
$ws.Cells.Item(95,5).Value = 754

# Row 98: Guayana Francesa -> Guayana Francesa
$ws.Cells.Item(98,2).Value = 9418
$ws.Cells.Item(98,3).Value = 31
$ws.Cells.Item(98,4).Value = 8986
$ws.Cells.Item(98,5).Value = 370

# Row 102: Gabon -> Gabon
$ws.Cells.Item(102,2).Value = 8621
$ws.Cells.Item(102,3).Value = 13
$ws.Cells.Item(102,4).Value = 7618
$ws.Cells.Item(102,5).Value = 950

# Row 118: Republica de Africa Central -> Mozambique
$ws.Cells.Item(118,1).Value = "Mozambique"
$ws.Cells.Item(118,2).Value = 4764
$ws.Cells.Item(118,3).Value = 117
$ws.Cells.Item(118,4).Value = 2763
$ws.Cells.Item(118,5).Value = 1973
$ws.Cells.Item(118,8).Value = 28

# Row 119: Mozambique -> Republica de Africa Central
$ws.Cells.Item(119,1).Value = "Republica de Africa Central"
$ws.Cells.Item(119,2).Value = 4735
$ws.Cells.Item(119,4).Value = 1825
$ws.Cells.Item(119,5).Value = 2848
$ws.Cells.Item(119,8).Value = 62

# Row 121: Ruanda -> Cuba
$ws.Cells.Item(121,1).Value = "Cuba"
$ws.Cells.Item(121,2).Value = 4459
$ws.Cells.Item(121,3).Value = 82
$ws.Cells.Item(121,4).Value = 3727
$ws.Cells.Item(121,5).Value = 628
$ws.Cells.Item(121,8).Value = 104

# Row 122: Surinam -> Ruanda
$ws.Cells.Item(122,1).Value = "Ruanda"
$ws.Cells.Item(122,2).Value = 4439
$ws.Cells.Item(122,4).Value = 2307
$ws.Cells.Item(122,5).Value = 2112
$ws.Cells.Item(122,8).Value = 20

# Row 123: Cuba -> Surinam
$ws.Cells.Item(123,1).Value = "Surinam"
$ws.Cells.Item(123,2).Value = 4419
$ws.Cells.Item(123,4).Value = 3595
$ws.Cells.Item(123,5).Value = 733
$ws.Cells.Item(123,8).Value = 91

# Row 128: Jamaica -> Siria
$ws.Cells.Item(128,1).Value = "Siria"
$ws.Cells.Item(128,2).Value = 3351
$ws.Cells.Item(128,3).Value = 62
$ws.Cells.Item(128,4).Value = 780
$ws.Cells.Item(128,5).Value = 2428
$ws.Cells.Item(128,7).Value = 3
$ws.Cells.Item(128,8).Value = 143

# Row 129: Eslovenia -> Jamaica
$ws.Cells.Item(129,1).Value = "Jamaica"
$ws.Cells.Item(129,2).Value = 3323
$ws.Cells.Item(129,3).Value = 140
$ws.Cells.Item(129,4).Value = 992
$ws.Cells.Item(129,5).Value = 2295
$ws.Cells.Item(129,7).Value = 2
$ws.Cells.Item(129,8).Value = 36

# Row 130: Gambia -> Eslovenia
$ws.Cells.Item(130,1).Value = "Eslovenia"
$ws.Cells.Item(130,2).Value = 3312
$ws.Cells.Item(130,3).Value = 79
$ws.Cells.Item(130,4).Value = 2587
$ws.Cells.Item(130,5).Value = 590
$ws.Cells.Item(130,8).Value = 135

# Row 131: Siria -> Gambia
$ws.Cells.Item(131,1).Value = "Gambia"
$ws.Cells.Item(131,2).Value = 3293
$ws.Cells.Item(131,3).Value = 18
$ws.Cells.Item(131,4).Value = 1460
$ws.Cells.Item(131,5).Value = 1734
$ws.Cells.Item(131,8).Value = 99

# Row 135: Mali -> Mali
$ws.Cells.Item(135,2).Value = 2898
$ws.Cells.Item(135,3).Value = 16
$ws.Cells.Item(135,4).Value = 2267
$ws.Cells.Item(135,5).Value = 503
$ws.Cells.Item(135,7).Value = 1
$ws.Cells.Item(135,8).Value = 128

# Row 155: Republica de Chipre -> Togo
$ws.Cells.Item(155,1).Value = "Togo"
$ws.Cells.Item(155,2).Value = 1528
$ws.Cells.Item(155,3).Value = 15
$ws.Cells.Item(155,4).Value = 1144
$ws.Cells.Item(155,5).Value = 348
$ws.Cells.Item(155,7).Value = 2
$ws.Cells.Item(155,8).Value = 36

# Row 156: Togo -> Republica de Chipre
$ws.Cells.Item(156,1).Value = "Republica de Chipre"
$ws.Cells.Item(156,2).Value = 1514
$ws.Cells.Item(156,3).Value = 3
$ws.Cells.Item(156,4).Value = 1237
$ws.Cells.Item(156,5).Value = 255
$ws.Cells.Item(156,8).Value = 22

# Row 157: Burkina Faso -> Burkina Faso
$ws.Cells.Item(157,2).Value = 1476
$ws.Cells.Item(157,3).Value = 10
$ws.Cells.Item(157,4).Value = 1118
$ws.Cells.Item(157,5).Value = 302

# Row 162: Principado de Andorra -> Principado de Andorra
$ws.Cells.Item(162,2).Value = 1301
$ws.Cells.Item(162,3).Value = 40
$ws.Cells.Item(162,4).Value = 938
$ws.Cells.Item(162,5).Value = 310

# Row 214: Montserrat -> Islas Malvinas
$ws.Cells.Item(214,1).Value = "Islas Malvinas"
$ws.Cells.Item(214,4).Value = 13
$ws.Cells.Item(214,8).Value = 0

# Row 215: Islas Malvinas -> Montserrat
$ws.Cells.Item(215,1).Value = "Montserrat"
$ws.Cells.Item(215,4).Value = 12
$ws.Cells.Item(215,8).Value = 1
